$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Hide the already-completed log rows (rows 2-13) ---------------------
$ws.Range("A2:A13").EntireRow.Hidden = $true

# --- 2. Add a new work-log entry in row 23 -----------------------------------
# First carry the previous "last row" formatting (row 22, before the edit)
# down onto the new last row (23)...
$ws.Range("A22:C22").Copy() | Out-Null
$ws.Range("A23:C23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ...then restyle the old last row (22) like a normal, non-terminal data row
# by copying the formatting from the row above it (21).
$ws.Range("A21:C21").Copy() | Out-Null
$ws.Range("A22:C22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Now fill in the new row's data.
$ws.Cells.Item(23, 1).Value = 45711
$ws.Cells.Item(23, 2).Value = 6
$ws.Cells.Item(23, 3).Value = "Final check with the tool and midterm report generation"

# --- 3. Update the view: drop the old scroll/selection, select the whole table
$ws.Range("A1:C23").Select() | Out-Null
